$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.33"
$ws.Range("E2").Value = "'-0.98%"
$ws.Range("D3").Value = "'27.13"
$ws.Range("E3").Value = "'3.15%"
$ws.Range("D4").Value = "'5.127"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("E5").Value = "'1.03%"
$ws.Range("D6").Value = "'6.475"
$ws.Range("E6").Value = "'-0.46%"
$ws.Range("D7").Value = "'0.8216"
$ws.Range("E7").Value = "'1.10%"
$ws.Range("D8").Value = "'0.8404"
$ws.Range("E8").Value = "'-0.06%"
$ws.Range("D9").Value = "'0.1328"
$ws.Range("E9").Value = "'-1.23%"
$ws.Range("D10").Value = "'0.06930"
$ws.Range("E10").Value = "'-0.30%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03142"
$ws.Range("E11").Value = "'0.79%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03018"
$ws.Range("E12").Value = "'7.42%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09396"
$ws.Range("E13").Value = "'0.17%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001514"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04218"
$ws.Range("E15").Value = "'-9.69%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005980"
$ws.Range("E16").Value = "'0.23%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006143"
$ws.Range("E17").Value = "'-0.67%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.515"
$ws.Range("E18").Value = "'-1.05%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.002"
$ws.Range("E19").Value = "'-1.33%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.225"
$ws.Range("E20").Value = "'5.06%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3113"
$ws.Range("E21").Value = "'-1.37%"
$ws.Range("D23").Value = "'3.566"
$ws.Range("E23").Value = "'-4.57%"
$ws.Range("E24").Value = "'-0.02%"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-2.14%"
$ws.Range("D26").Value = "'0.004463"
$ws.Range("E26").Value = "'-3.29%"
$ws.Range("D27").Value = "'0.00009802"
$ws.Range("E27").Value = "'2.19%"
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("E28").Value = "'39.53%"
$ws.Range("D40").Value = "'0.03653"
$ws.Range("E40").Value = "'-0.16%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1054"
$ws.Range("E41").Value = "'-22.25%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002530"
$ws.Range("E42").Value = "'-4.79%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003437"
$ws.Range("E43").Value = "'-44.17%"
$ws.Range("D44").Value = "'0.008974"
$ws.Range("E44").Value = "'-0.15%"
$ws.Range("D45").Value = "'0.00005312"
$ws.Range("E45").Value = "'0.46%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("D48").Value = "'0.002643"
$ws.Range("E48").Value = "'28.28%"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E50").Value = "'0.06%"